# Update the "Förändrad" (Changed) date column (C) from 2026-02-22 (serial 46075)
# to 2026-02-23 (serial 46076) for every data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2() -eq 46075) {
        $cell.Value = 46076
    }
}
